$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The formulas in B2:F8 are being replaced by their plain computed values
# (i.e. converted from formulas to static numbers), cell by cell so the
# COM layer treats each as a scalar rather than an array.
for ($r = 2; $r -le 8; $r++) {
    for ($c = 2; $c -le 6; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value2
        $cell.Formula = ""
        $cell.Value2 = $v
    }
}

# A few values were also corrected by the author when the formulas were
# dropped; apply those corrected numbers explicitly.
$ws.Cells.Item(2, 6).Value2 = 1804978633        # F2
$ws.Cells.Item(7, 2).Value2 = 11176.5           # B7
$ws.Cells.Item(7, 3).Value2 = 13027.3           # C7

# Update the active selection to match the saved state
$ws.Range("E9").Select()
